# Documentatie - Timeline.xlsx edit
# "made bookmark for 4 and 5 + finished 4.1"
#
# - Adds a new timeline entry on row 14 (day 7, 16/6/2024, 3.33 hours,
#   "Creare curpins pentru capitolele 4 si 5 + Finalizare 4.1")
# - Inserts 3 new blank rows before the totals block, pushing it from
#   rows 16-18 down to rows 19-21
# - Extends the "Total hours Spent" SUM formula to cover the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the totals block (old rows 16:18),
# pushing the totals rows down to 19:21 (mirrors what Excel does, also
# shifting the mergeCells C16:C18/D16:D18 to C19:C21/D19:D21).
$ws.Rows("16:18").Insert()

# Fill in the new timeline entry on (still) row 14.
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "16/6/2024"
$ws.Range("C14").Value = 3.33
$ws.Range("D14").Value = "Creare curpins pentru capitolele 4 si 5 + Finalizare 4.1"

# Extend the total-hours formula (now on row 19) to include the new data
# rows up through row 18.
$ws.Range("D19").Formula = "=SUM(C3:C18)"

# Re-fit the row height back to the sheet default (setting the formula
# above can otherwise leave an explicit custom row height behind).
$ws.Rows(19).AutoFit()

# Match the saved selection/active cell.
$ws.Range("D18").Select()
